# concerts.xlsx report update: extend the "Три дня дождя" (tour "Волшебство") block with
# three more concert rows and rebuild the "Kai Angel" block further down, adding a new
# track/concert and leaving row 7 blank as the section separator.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 held the old "Название тура" / "Кристаллические лярвы" header for the second
# artist block; that block now starts two rows lower (8-9), so clear the stale cells
# and leave row 7 as the blank separator row.
$ws.Range("A7:B7").ClearContents()

# Row 8 used to be a full 6-column concert row ("Название концерта"..."Вача"); in the
# new layout it is just the 2-cell "Исполнитель" header, so drop the leftover C:F cells.
$ws.Range("C8:F8").ClearContents()

$ws.Range("A1").Value = "Исполнитель"
$ws.Range("B1").Value = "Три дня дождя"

$ws.Range("A2").Value = "Название тура"
$ws.Range("B2").Value = "Волшебство"

$ws.Range("A3").Value = "Название концерта"
$ws.Range("B3").Value = "Театр имени Достоевского"
$ws.Range("C3").Value = "Дата"
$ws.Range("D3").Value = "2024-05-31 21:00"
$ws.Range("E3").Value = "Город"
$ws.Range("F3").Value = "Муром"

$ws.Range("A4").Value = "Название концерта"
$ws.Range("B4").Value = "Milo Hall"
$ws.Range("C4").Value = "Дата"
$ws.Range("D4").Value = "2024-06-08 21:01"
$ws.Range("E4").Value = "Город"
$ws.Range("F4").Value = "Кулебаки"

$ws.Range("A5").Value = "Название концерта"
$ws.Range("B5").Value = "Supe Hall"
$ws.Range("C5").Value = "Дата"
$ws.Range("D5").Value = "2024-05-30 21:01"
$ws.Range("E5").Value = "Город"
$ws.Range("F5").Value = "Вача"

$ws.Range("A6").Value = "Название концерта"
$ws.Range("B6").Value = "Концертный зал Мило"
$ws.Range("C6").Value = "Дата"
$ws.Range("D6").Value = "2024-06-04 21:09"
$ws.Range("E6").Value = "Город"
$ws.Range("F6").Value = "Нижний Новгород"

$ws.Range("A8").Value = "Исполнитель"
$ws.Range("B8").Value = "Kai Angel"

$ws.Range("A9").Value = "Название тура"
$ws.Range("B9").Value = "Russian Underground"

$ws.Range("A10").Value = "Название концерта"
$ws.Range("B10").Value = "Supe Hall"
$ws.Range("C10").Value = "Дата"
$ws.Range("D10").Value = "2024-05-31 21:02"
$ws.Range("E10").Value = "Город"
$ws.Range("F10").Value = "Вача"

$ws.Range("A11").Value = "Название концерта"
$ws.Range("B11").Value = "Концертный зал Мило"
$ws.Range("C11").Value = "Дата"
$ws.Range("D11").Value = "2024-06-02 21:02"
$ws.Range("E11").Value = "Город"
$ws.Range("F11").Value = "Нижний Новгород"

